$wb = $excel.ActiveWorkbook

# --- 1) Rename sheet "Cartes" -> "Feuille" ---
$wsCartes = $wb.Worksheets.Item("Cartes")
$wsCartes.Name = "Feuille"

# --- 1b) Drop a stray leftover empty cell on "09_03_2025" ---
$ws3 = $wb.Worksheets.Item("09_03_2025")
$ws3.Range("G19").Value = ""

# --- 2) Apply data corrections to the "10_03_2025" sheet ---
$ws = $wb.Worksheets.Item("10_03_2025")

# Fix apostrophe-containing names to use the literal "&apos;" text
# (matching the convention already used on the other sheets)
$ws.Range("A3").Value = "Reverse World&apos;s Giratina Lv.60"
$ws.Range("A4").Value = "Icy Sky&apos;s Shaymin Lv.62"
$ws.Range("A6").Value = "Alto Mare&apos;s Latias"
$ws.Range("A7").Value = "Alto Mare&apos;s Latios"

# Clear stray empty cells in column B
$ws.Range("B6").Value = ""
$ws.Range("B7").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("B19").Value = ""
$ws.Range("B20").Value = ""
$ws.Range("B21").Value = ""

# Corrected average-price values in column G
$ws.Range("G2").Value = 227.62
$ws.Range("G4").Value = 6.95
$ws.Range("G6").Value = 26.63
$ws.Range("G9").Value = 3.76
$ws.Range("G14").Value = 27.52
$ws.Range("G16").Value = 8.12
$ws.Range("G18").Value = 14.99

# --- 3) Duplicate the corrected "10_03_2025" sheet into a new "11_03_2025" sheet ---
$ws.Copy($null, $ws)
$wsNew = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew.Name = "11_03_2025"
